# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" sheet right after "总计" (before "2022-Q3"),
#    carrying the per-fund holdings table for that quarter.
# 2) Insert a new row into "总计" (the summary sheet) for 2022-Q4,
#    above the existing 2022-Q3 row, shifting everything else down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: "总计" summary sheet - insert the 2022-Q4 row at the top of
# the data (row 2), pushing the other quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

$summary.Range("A2").Value = 0
$summary.Range("A2").Font.Bold = $true
$summary.Range("A2").Borders.LineStyle = 1
$summary.Range("A2").HorizontalAlignment = -4108
$summary.Range("A2").VerticalAlignment = -4160

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.33

# The row-insert shifted the other quarters down but kept their old
# "row index" values (column A) unchanged (0,1,2,3,4) - bump each by
# one so the sequential numbering stays consistent (1,2,3,4,5).
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

# ---------------------------------------------------------------------
# Part 2: brand-new "2022-Q4" worksheet with the fund holdings detail,
# placed immediately after "总计".
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").Borders.LineStyle = 1
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160

# The "基金代码" (fund code, to keep leading zeros) and the
# "基金规模"/"股票总仓位"/"仓位占比"/"持有市值(亿元)" columns all hold
# text (not numeric) values in this workbook, so force the number
# format to Text before writing them to avoid COM auto-converting the
# numeric-looking strings into real numbers.
$q4.Range("B2:B5").NumberFormat = "@"
$q4.Range("D2:G5").NumberFormat = "@"

$rows = @(
    @(0, "010852", "中欧内需成长混合A",   "3.73", "93.04", "4.90", "0.1828", 9),
    @(1, "162203", "泰达宏利稳定混合",     "2.98", "91.72", "3.55", "0.1058", 8),
    @(2, "010853", "中欧内需成长混合C",   "0.54", "93.04", "4.90", "0.0265", 9),
    @(3, "620004", "金元顺安价值增长混合", "0.52", "77.68", "1.97", "0.0102", 4)
)

$r = 2
foreach ($row in $rows) {
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("A$r").Font.Bold = $true
    $q4.Range("A$r").Borders.LineStyle = 1
    $q4.Range("A$r").HorizontalAlignment = -4108
    $q4.Range("A$r").VerticalAlignment = -4160

    $q4.Range("B$r").Value = $row[1]
    $q4.Range("C$r").Value = $row[2]
    $q4.Range("D$r").Value = $row[3]
    $q4.Range("E$r").Value = $row[4]
    $q4.Range("F$r").Value = $row[5]
    $q4.Range("G$r").Value = $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}
